$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlLeft = -4131
$xlCenter = -4108

# ---------------------------------------------------------------------------
# Helper: write a value into a cell and apply the standard left/center
# alignment used throughout this sheet's data rows.
# ---------------------------------------------------------------------------
function Set-IssueCell($Range, $Value, $Wrap) {
    if ($null -ne $Value) {
        $Range.Value = $Value
    }
    $Range.HorizontalAlignment = $xlLeft
    $Range.VerticalAlignment = $xlCenter
    $Range.WrapText = $Wrap
}

# Helper that forces a "looks like a number" string (e.g. "0.8") to be
# stored as text instead of being auto-converted to a numeric value.
function Set-IssueTextCell($Range, $Text, $Wrap) {
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = "Normal"
    $Range.HorizontalAlignment = $xlLeft
    $Range.VerticalAlignment = $xlCenter
    $Range.WrapText = $Wrap
}

# ---------------------------------------------------------------------------
# Seed the new shared strings in the same order they were added originally:
# the two new issue descriptions first, then the new version label "0.8".
# ---------------------------------------------------------------------------
Set-IssueCell $ws.Range("B11") "Pantalla 5: Pantalla de cabecera: La versión situarla arriba de la primer letra del nombre del usuario conectado." $true
Set-IssueCell $ws.Range("B12") "Pantalla 5: Campo de Moneda. Entendemos que es una lista desplegable con las monedas configuradas en Netsuite." $true

# ---------------------------------------------------------------------------
# Row 10: the version value "0.7" becomes the new text value "0.8"
# ---------------------------------------------------------------------------
Set-IssueTextCell $ws.Range("D10") "0.8" $false

# ---------------------------------------------------------------------------
# Row 11 (new) - new issue about the version placement on screen 5
# ---------------------------------------------------------------------------
Set-IssueCell $ws.Range("A11") 10 $false
Set-IssueCell $ws.Range("C11") "PUC" $false
Set-IssueTextCell $ws.Range("D11") "0.8" $false
Set-IssueCell $ws.Range("E11") "Pendiente aceptación" $true
Set-IssueCell $ws.Range("F11") $null $false
$ws.Rows.Item(11).RowHeight = 30

# ---------------------------------------------------------------------------
# Row 12 (new) - new issue about the currency field
# ---------------------------------------------------------------------------
Set-IssueCell $ws.Range("A12") 11 $false
Set-IssueCell $ws.Range("C12") "PUC" $false
Set-IssueTextCell $ws.Range("D12") "0.8" $false
Set-IssueCell $ws.Range("E12") "Pendiente aceptación" $true
Set-IssueCell $ws.Range("F12") $null $false
$ws.Rows.Item(12).RowHeight = 45

# ---------------------------------------------------------------------------
# Grow the worksheet table ("Tabla1") so it covers the two new rows
# ---------------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:F12"))

# ---------------------------------------------------------------------------
# Match the saved selection state
# ---------------------------------------------------------------------------
$ws.Range("C12").Select()
